# Update "想去人数" (F column) counts for sheets "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 25
    6  = 93
    7  = 467
    9  = 25
    10 = 584
    11 = 33
    12 = 310
    14 = 380
    16 = 97
    17 = 13
    21 = 976
    22 = 1411
    23 = 307
    28 = 44
    29 = 88
    30 = 228
    31 = 259
    33 = 1634
    37 = 589
    38 = 298
    39 = 3737
    40 = 435
    41 = 208
    42 = 923
    45 = 76
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
